$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for dates 2021-04-18 through 2021-04-21 (serials 44304-44307)
$newRows = @(
    @{ Row = 230; A = 44304; B = 7; C = 52; D = 203.5304708599162 },
    @{ Row = 231; A = 44305; B = 4; C = 49; D = 191.7883283103057 },
    @{ Row = 232; A = 44306; B = 8; C = 51; D = 199.6164233433794 },
    @{ Row = 233; A = 44307; B = 4; C = 52; D = 203.5304708599162 }
)

# Copy the formatting of the last existing date cell (A229, style s="2")
# so the new date cells pick up the same style instead of minting a new one.
$ws.Range("A229").Copy()

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).PasteSpecial(-4122)
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
}
